# Applies an Ofsted CSC ILACS data refresh to the
# "ofsted_cs_inspections_overview" sheet:
#   - Warwickshire (la_code 937) received a brand-new inspection: new
#     report link/dates, new grades (good -> requires improvement across
#     most columns), a new (unnamed) inspector, and 0 recorded inspections
#     for that inspector so far.
#   - Every other local authority that lists Warwickshire (937) as a
#     statistical neighbour has its "stat_neighbour_judgement" text
#     updated so the (937, 'good') tuple now reads (937, 'requires
#     improvement').
#   - Alison Smale, who used to be credited with Warwickshire's inspection,
#     now has one fewer inspection to her name (6 -> 5) everywhere she is
#     still listed as the inspector.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes $text into $cell while keeping it a genuine text cell: a plain
# `.Value = $text` assignment lets Excel's COM layer auto-convert
# date-/number-looking strings (e.g. "02/06/2025", "5") into real
# dates/numbers. Forcing the number format to Text ("@") first suppresses
# that, then ClearFormats() drops the now-unneeded explicit format so the
# cell's style stays exactly as it was (no stray numFmt / style index).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# --- locate columns by header text (row 1) so this is resilient to column
#     reordering -------------------------------------------------------
$lastCol = $ws.Cells.Item(1, $ws.Columns.Count).End(-4159).Column
$col = @{}
for ($c = 1; $c -le $lastCol; $c++) {
    $header = [string]$ws.Cells.Item(1, $c).Value2
    $col[$header] = $c
}

$colNeighbours   = $col["stat_neighbours"]
$colNeighJudge   = $col["stat_neighbour_judgement"]
$colAuthority    = $col["local_authority"]
$colLink         = $col["inspection_link"]
$colOverall      = $col["overall_effectiveness_grade"]
$colFramework    = $col["inspection_framework"]
$colInspector    = $col["inspector_name"]
$colStartDate    = $col["inspection_start_date"]
$colEndDate      = $col["inspection_end_date"]
$colPubDate      = $col["publication_date"]
$colImpact       = $col["impact_of_leaders_grade"]
$colHelp         = $col["help_and_protection_grade"]
$colInCare       = $col["in_care_grade"]
$colLeavers      = $col["care_leavers_grade"]
$colInspCount    = $col["inspectors_inspections_count"]

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# --- 1) Warwickshire's own row: new inspection details --------------
$warwickshireRow = $ws.Cells.Find("warwickshire", [System.Type]::Missing, [System.Type]::Missing, 1).Row

$ws.Cells.Item($warwickshireRow, $colLink).Value = "https://files.ofsted.gov.uk/v1/file/50289041"
Set-TextValue $ws.Cells.Item($warwickshireRow, $colStartDate) "02/06/2025"
Set-TextValue $ws.Cells.Item($warwickshireRow, $colEndDate)   "06/06/2025"
Set-TextValue $ws.Cells.Item($warwickshireRow, $colPubDate)   "14/10/25"

$ws.Cells.Item($warwickshireRow, $colOverall).Value   = "requires improvement"
$ws.Cells.Item($warwickshireRow, $colFramework).Value = "short"
$ws.Cells.Item($warwickshireRow, $colInspector).Value = "None"
$ws.Cells.Item($warwickshireRow, $colImpact).Value    = "requires improvement"
$ws.Cells.Item($warwickshireRow, $colHelp).Value      = "requires improvement"
$ws.Cells.Item($warwickshireRow, $colInCare).Value    = "requires improvement"
$ws.Cells.Item($warwickshireRow, $colLeavers).Value   = "good"
Set-TextValue $ws.Cells.Item($warwickshireRow, $colInspCount) "0"

# --- 2) Every LA whose statistical-neighbour list includes 937 (Warwickshire)
#        gets its neighbour-judgement text updated in place -----------
for ($r = 2; $r -le $lastRow; $r++) {
    $neighbours = [string]$ws.Cells.Item($r, $colNeighbours).Value2
    if ([string]::IsNullOrEmpty($neighbours)) { continue }

    $isNeighbour = $false
    foreach ($code in ($neighbours -split ",")) {
        if ($code.Trim() -eq "937") { $isNeighbour = $true }
    }

    if ($isNeighbour) {
        $judgeCell = $ws.Cells.Item($r, $colNeighJudge)
        $judgeText = [string]$judgeCell.Value2
        $updated = $judgeText.Replace("(937, 'good')", "(937, 'requires improvement')")
        if ($updated -ne $judgeText) {
            $judgeCell.Value = $updated
        }
    }
}

# --- 3) Alison Smale loses Warwickshire's inspection from her tally ---
for ($r = 2; $r -le $lastRow; $r++) {
    if ($r -eq $warwickshireRow) { continue }
    $inspector = [string]$ws.Cells.Item($r, $colInspector).Value2
    if ($inspector -eq "alison smale") {
        Set-TextValue $ws.Cells.Item($r, $colInspCount) "5"
    }
}
